$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Database" category label to "DB" for rows 12-16 (column B)
$ws.Range("B12:B16").Value = "DB"

# Update the text of A15 from "Database migration..." to "DB migration..."
$ws.Range("A15").Value = "DB migration script failed due to syntax error."

# Restore the selection to G10 as recorded in the sheet view
$ws.Range("G10").Select()
